$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View state (best-effort; window/scroll chrome is largely host-only state) ---
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 4650
$ws.Range("F24").Select() | Out-Null

# --- Four volunteers newly marked as signed up (E column = 1) with a mailto
#     hyperlink added to their e-mail address in column C ---

# Row 24 - Emmelie Simoens
$ws.Hyperlinks.Add($ws.Range("C62"), "mailto:Sarah.Slabbaert@UGent.be") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C58"), "mailto:Rebecca.Willems@UGent.be") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C27"), "mailto:Frederik.DeSpiegeleer@UGent.be") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C24"), "mailto:Emmelie.Simoens@UGent.be") | Out-Null

# Re-apply the existing "hyperlink" cell style from an already-linked row (C8)
# so the style index matches the rest of the sheet instead of minting a new one.
$ws.Range("C8").Copy()
$ws.Range("C62").PasteSpecial(-4122)
$ws.Range("C58").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C24").PasteSpecial(-4122)

# Give the four E cells the same numeric style as the existing filled-in E column
# cells (E8, numFmtId 2) and set their value to 1 (signed up).
$ws.Range("E8").Copy()
$ws.Range("E62").PasteSpecial(-4122)
$ws.Range("E58").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("E24").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("E58").Value = 1
$ws.Range("E62").Value = 1
